# Swap the order of "dnasr281@gmail.com" and "System" in the
# "Recorded By" column (G) wherever both names are combined together.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
